$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-15 19:48:36'
$ws.Range("E3").Value = '2026-02-15 19:48:38'
$ws.Range("I3").Value = '1.3 mm'
$ws.Range("O3").Value = '-5.2 °C'
$ws.Range("E4").Value = '2026-02-15 19:48:40'
$ws.Range("E5").Value = '2026-02-15 19:48:43'
$ws.Range("I5").Value = '4.8 mm'
$ws.Range("O5").Value = '-4.6 °C'
$ws.Range("E6").Value = '2026-02-15 19:48:45'
$ws.Range("E7").Value = '2026-02-15 19:48:48'
$ws.Range("O7").Value = '11.7 °C'
$ws.Range("E8").Value = '2026-02-15 19:48:50'
$ws.Range("E9").Value = '2026-02-15 19:48:53'
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = '50%'
$ws.Range("E10").Value = '2026-02-15 19:48:55'
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = '70%'
$ws.Range("E11").Value = '2026-02-15 19:48:58'
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = '42%'
$ws.Range("O11").Value = '7.4 °C'
$ws.Range("E12").Value = '2026-02-15 19:49:00'
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = '55%'
$ws.Range("E13").Value = '2026-02-15 19:49:02'
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = '35%'
$ws.Range("J13").Value = '1015.4 hPa'
$ws.Range("E14").Value = '2026-02-15 19:49:05'
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = '59%'
$ws.Range("E15").Value = '2026-02-15 19:49:07'
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = '50%'
$ws.Range("O15").Value = '10.7 °C'
$ws.Range("E16").Value = '2026-02-15 19:49:10'
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = '61%'
$ws.Range("O16").Value = '-2.0 °C'
$ws.Range("E17").Value = '2026-02-15 19:49:12'
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = '37%'
$ws.Range("E18").Value = '2026-02-15 19:49:14'
$ws.Range("E19").Value = '2026-02-15 19:49:17'
$ws.Range("E20").Value = '2026-02-15 19:49:19'
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = '60%'
$ws.Range("E21").Value = '2026-02-15 19:49:21'
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = '38%'
$ws.Range("E22").Value = '2026-02-15 19:49:24'
$ws.Range("N22").Value = '-6.5 °C 19:06 TU'
$ws.Range("E23").Value = '2026-02-15 19:49:26'
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = '65%'
$ws.Range("I23").Value = '2.2 mm'
$ws.Range("E24").Value = '2026-02-15 19:49:29'
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = '67%'
$ws.Range("E25").Value = '2026-02-15 19:49:31'
$ws.Range("M25").Value = '2.6 °C 19:03 TU'
$ws.Range("O25").Value = '-1.6 °C'
$ws.Range("E26").Value = '2026-02-15 19:49:34'
$ws.Range("E27").Value = '2026-02-15 19:49:36'
$ws.Range("E28").Value = '2026-02-15 19:49:38'
$ws.Range("E29").Value = '2026-02-15 19:49:41'
$ws.Range("E30").Value = '2026-02-15 19:49:43'
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = '54%'
$ws.Range("E31").Value = '2026-02-15 19:49:46'
$ws.Range("E32").Value = '2026-02-15 19:49:48'
$ws.Range("O32").Value = '3.7 °C'
$ws.Range("E33").Value = '2026-02-15 19:49:50'
$ws.Range("E34").Value = '2026-02-15 19:49:53'
$ws.Range("M34").Value = '4.6 °C 19:12 TU'
$ws.Range("O34").Value = '1.0 °C'
$ws.Range("E35").Value = '2026-02-15 19:49:55'
$ws.Range("O35").Value = '4.1 °C'
$ws.Range("E36").Value = '2026-02-15 19:49:58'
$ws.Range("E37").Value = '2026-02-15 19:50:00'
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = '53%'
$ws.Range("J37").Value = '1016.3 hPa'
$ws.Range("E38").Value = '2026-02-15 19:50:02'
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H38").Value = '67%'
$ws.Range("E39").Value = '2026-02-15 19:50:05'
$ws.Range("O39").Value = '-2.9 °C'
$ws.Range("E40").Value = '2026-02-15 19:50:07'
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = '36%'
$ws.Range("O40").Value = '8.9 °C'
$ws.Range("E41").Value = '2026-02-15 19:50:10'
$ws.Range("E42").Value = '2026-02-15 19:50:12'
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = '56%'
$ws.Range("O42").Value = '10.6 °C'
$ws.Range("E43").Value = '2026-02-15 19:50:14'
$ws.Range("K43").Value = '12.8 MJ/m2'
$ws.Range("O43").Value = '6.3 °C'
$ws.Range("E44").Value = '2026-02-15 19:50:17'
$ws.Range("H44").NumberFormat = "@"
$ws.Range("H44").Value = '77%'
$ws.Range("I44").Value = '2.3 mm'
$ws.Range("O44").Value = '-4.1 °C'
$ws.Range("E45").Value = '2026-02-15 19:50:19'
$ws.Range("I45").Value = '1.0 mm'
$ws.Range("J45").Value = '1023.4 hPa'
$ws.Range("O45").Value = '1.0 °C'
$ws.Range("E46").Value = '2026-02-15 19:50:21'
